$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (new weekly entry replacing the previous top row) ---
$ws.Range("D2").Value = 44425
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("P2").Value = 560

# --- Update row 3 (shifted down from what used to be row 2) ---
$ws.Range("D3").Value = 44340
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("P3").Value = 600

# --- Row 4 (new entry) ---
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44421
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100112026
$ws.Range("G4").Value = "Haba"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 15000
$ws.Range("N4").Value = "$/saco 25 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 600
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"

# --- Row 5 (entry that used to be row 3) ---
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value = 44376
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 100112026
$ws.Range("G5").Value = "Haba"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 15
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = "$/saco 25 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 480
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"

# --- Row 6 (new entry) ---
$ws.Range("A6").Value = 12
$ws.Range("B6").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C6").Value = "Metropolitana"
$ws.Range("D6").Value = 44418
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 12
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 600
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
